$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: "ImmSrc" -> "ImmType" ---
$ws.Range("F2").Value = "ImmType"

# --- F3: rich text change ---
# Old: "1 [" (regular) + "addi should sign extend immediate" (italic) + "]" (regular)
# New: "000 [" (regular) + "I-type sign extend 12 (instruction[31-20]) bits to 32 bits" (italic) + "]" (regular)
$cell = $ws.Range("F3")
$newText = "000 [I-type sign extend 12 (instruction[31-20]) bits to 32 bits]"
$cell.Value = $newText
$cell.Characters(6, 58).Font.Italic = $true
$cell.Characters(64, 1).Font.Italic = $false

# --- Column F width: widen to fit the new, longer text ---
$ws.Columns("F").ColumnWidth = 50.8

# --- Update the selected/active cell in the bottom-right frozen pane ---
[void]$ws.Range("E14").Select()
